# Update the "Channel Tests" sheet with new channel test cases and remove
# the old abstract event test references (per commit message), mirroring
# the target OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Channel Tests")

# Fill the Test Id column for all new rows up front.
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8
$ws.Range("A10").Value = 9
$ws.Range("A11").Value = 10
$ws.Range("A12").Value = 11
$ws.Range("A13").Value = 12
$ws.Range("A14").Value = 13
$ws.Range("A15").Value = 14
$ws.Range("A16").Value = 15
$ws.Range("A17").Value = 16
$ws.Range("A18").Value = 17
$ws.Range("A19").Value = 18

# Descriptions / Expected / Actual results, written in the same order the
# shared-string table grows in the target workbook.
$ws.Range("B5").Value = "Create new channel with default system controller, valid id, valid sensor name and not active"
$ws.Range("B6").Value = "Create new channel with default system controller, valid id, valid sensor name and is active"
$ws.Range("B8").Value = 'Add sensor to the channel with string of "GATE"'
$ws.Range("B9").Value = 'Add sensor to channel with string of "bad sensor"'
$ws.Range("C9").Value = "Illegal argument Exception for trying to add an invalid sensor type"
$ws.Range("C8").Value = "Channel adds a sensor of type gate"
$ws.Range("B7").Value = "Add a sensor to the channel with a sensor of type PAD"
$ws.Range("C7").Value = "Channel's sensor is set to a PAD"

$ws.Range("D2").Value = "Same as expected"
$ws.Range("D3").Value = "Same as expected"
$ws.Range("D4").Value = "Same as expected"

$ws.Range("D5").Value = "Same as expedted"
$ws.Range("C5").Value = "Channel is created and is not active(false)"
$ws.Range("C6").Value = "Channel is created and is active(true)"

$ws.Range("D6").Value = "same as expected"
$ws.Range("D7").Value = "same as expected"
$ws.Range("D8").Value = "same as expected"
$ws.Range("D9").Value = "same as expected"

$ws.Range("B10").Value = "A channel with sensor of type PAD has sensor disconnected"
$ws.Range("C10").Value = "Channel's sensor is null and sensor type is null"
$ws.Range("D10").Value = "same as expected"

$ws.Range("B11").Value = "Trigger channel's sensor that is not active(false)"
$ws.Range("C11").Value = "The sensor is not triggered because channel is not active"
$ws.Range("D11").Value = "same as expected"

$ws.Range("B12").Value = "Trigger a channel's sensor that is active(true) but sensor is null"
$ws.Range("C12").Value = "the sensor is not triggered because it does not exist"
$ws.Range("D12").Value = "same as expected"

$ws.Range("B13").Value = "Test get channel state when not active"
$ws.Range("C13").Value = "The state is not active so false"
$ws.Range("D13").Value = "same as expected"

$ws.Range("B14").Value = "test get channel state when active"
$ws.Range("C14").Value = "The state is true because it is active"
$ws.Range("D14").Value = "same as expected"

$ws.Range("B15").Value = "Get channel id with id 0f 100"
$ws.Range("C15").Value = "id is returned to be 100"
$ws.Range("D15").Value = "same as expected"

$ws.Range("B16").Value = "Channel exit"
$ws.Range("C16").Value = "Id is -1, active is false and current sensor is set to null if there was on"
$ws.Range("D16").Value = "same as expected"

$ws.Range("B17").Value = "Channel trigger sensor no sensor"
$ws.Range("C17").Value = "cannot trigger sensor since nul"
$ws.Range("D17").Value = "same as expected"

# Re-fit the description (B) and new actual-results (D) columns to their
# widened content; column C's best-fit width is unchanged from the source
# workbook so it is left untouched.
$ws.Columns.Item(2).ColumnWidth = 76
$ws.Columns.Item(4).ColumnWidth = 15

# Scroll/selection state mirrors the saved view in the target workbook.
$ws.Range("B18").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 2
